$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Rename "Act Q1" / "Act Q2" / "Act Q3" -> "Act Q4" / "Act Q5" / "Act Q6"
#    on slides 7, 8, 9 (only the "Qn" portion is retyped, leaving the
#    "Act " run intact, same as selecting just the digits in the UI).
# ------------------------------------------------------------------
$actRenames = @{ 7 = "Q4"; 8 = "Q5"; 9 = "Q6" }
foreach ($idx in $actRenames.Keys) {
    $slide = $p.Slides.Item($idx)
    $titleTr = $slide.Shapes.Item(1).TextFrame.TextRange
    $len = $titleTr.Text.Length
    $sub = $titleTr.Characters($len - 1, 2)
    $sub.Text = $actRenames[$idx]
}

# ------------------------------------------------------------------
# 2) Add three new "Title and Content" slides (10, 11, 12) at the end.
# ------------------------------------------------------------------

# --- Slide 10: Assert Q7 ---
$s10 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Assert Q7"

$body10 = $s10.Shapes.Item(2).TextFrame.TextRange
$body10.Text = "Take your answer from Q6 and let’s check if it returns true:`rTip:`rSince we are returning a “TaskExecutionResult” we need to assert that type, not a bool."
$body10.Paragraphs(2, 1).IndentLevel = 2
$body10.Paragraphs(3, 1).IndentLevel = 3
$body10.Paragraphs(1, 1).Font.Size = 28
$body10.Paragraphs(2, 1).Font.Size = 24
$body10.Paragraphs(3, 1).Font.Size = 24

# --- Slide 11: Assert Q8 ---
$s11 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "Assert Q8"

$body11 = $s11.Shapes.Item(2).TextFrame.TextRange
$line1 = "What if we have a collection that we want to compare?"
$line2 = "Using “Assert.That”, compare IEnumerable answers to IEnumerable results. "
$body11.Text = "$line1`r$line2"
$body11.Paragraphs(2, 1).IndentLevel = 2
$body11.Paragraphs(1, 1).Font.Size = 28
$body11.Paragraphs(2, 1).Font.Size = 24

# Underline + shadow on "answers" within the second paragraph.
$answersStart = $line1.Length + 1 + $line2.IndexOf("answers") + 1
$answersRange = $body11.Characters($answersStart, "answers".Length)
$answersRange.Font.Underline = $true
$answersRange.Font.Shadow = $true

# Underline on "results" within the second paragraph.
$resultsStart = $line1.Length + 1 + $line2.IndexOf("results") + 1
$resultsRange = $body11.Characters($resultsStart, "results".Length)
$resultsRange.Font.Underline = $true

# --- Slide 12: Assert Q9 ---
$s12 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = "Assert Q9"

$body12 = $s12.Shapes.Item(2).TextFrame.TextRange
$full12 = "Take the bool isCorrect and assert that it’s false."
$body12.Text = $full12
$body12.Font.Size = 28
$body12.TextFrame.AutoSize = 2

$isCorrectStart = $full12.IndexOf("isCorrect") + 1
$isCorrectRange = $body12.Characters($isCorrectStart, "isCorrect".Length)
$isCorrectRange.Font.Underline = $true
